{"js": "// Replace each paragraph's text in document order according to the\n// mapping captured from the authoritative diff. The document body is a\n// flat sequence of 101 paragraphs in this file: 1 title/date paragraph\n// followed by 100 arithmetic-problem paragraphs (one per table cell).\n// Both old and new text are recorded per-position so we can verify the\n// current text before mutating it (defensive against drift) and still\n// fall back to positional replacement if matching fails.\nconst REPLACEMENTS = [[\"2026-01-06 Tuesday\", \"2026-01-07 Wednesday\"], [\"77+6=\", \"69-54=\"], [\"46+47=\", \"40+21=\"], [\"79-0=\", \"25+60=\"], [\"74-36=\", \"61-38=\"], [\"85-84=\", \"29+50=\"], [\"5+52=\", \"77-31=\"], [\"33+8=\", \"25+37=\"], [\"50-10=\", \"3+63=\"], [\"20+65=\", \"71-1=\"], [\"76-29=\", \"52-12=\"], [\"96-45=\", \"89-70=\"], [\"55-3=\", \"54-15=\"], [\"47-37=\", \"33+33=\"], [\"23+5=\", \"21+11=\"], [\"53+5=\", \"41+2=\"], [\"47+2=\", \"76-50=\"], [\"48-10=\", \"18+29=\"], [\"79-51=\", \"9+55=\"], [\"53-42=\", \"32+25=\"], [\"21+15=\", \"59+33=\"], [\"7+60=\", \"49+30=\"], [\"46+11=\", \"42+26=\"], [\"93-77=\", \"51+0=\"], [\"31+52=\", \"40+11=\"], [\"84-1=\", \"53-0=\"], [\"23+60=\", \"96-7=\"], [\"40-4=\", \"42-30=\"], [\"90-63=\", \"86-1=\"], [\"35+32=\", \"72+1=\"], [\"50-20=\", \"73-32=\"], [\"49+7=\", \"88-32=\"], [\"87-66=\", \"68+12=\"], [\"88-51=\", \"82+3=\"], [\"64+20=\", \"6+51=\"], [\"13+68=\", \"70+14=\"], [\"5+46=\", \"93-87=\"], [\"15+38=\", \"53-8=\"], [\"15+83=\", \"28-7=\"], [\"72+21=\", \"68-31=\"], [\"66-18=\", \"2-0=\"], [\"29-0=\", \"14+25=\"], [\"91-64=\", \"89+2=\"], [\"42+8=\", \"85-42=\"], [\"22+56=\", \"53-25=\"], [\"41-9=\", \"66+27=\"], [\"21-18=\", \"31+19=\"], [\"52+35=\", \"33+25=\"], [\"46-7=\", \"30+39=\"], [\"85-21=\", \"47-6=\"], [\"54-5=\", \"47+7=\"], [\"93+3=\", \"8+57=\"], [\"53+8=\", \"20-12=\"], [\"66-30=\", \"29+19=\"], [\"96-59=\", \"72-18=\"], [\"81-41=\", \"38+7=\"], [\"99-37=\", \"77-6=\"], [\"62-40=\", \"18+4=\"], [\"98-6=\", \"56-51=\"], [\"53-1=\", \"33+15=\"], [\"95-63=\", \"30-11=\"], [\"11+65=\", \"17+50=\"], [\"19+21=\", \"16-15=\"], [\"46+43=\", \"19+8=\"], [\"66-19=\", \"19-19=\"], [\"8+59=\", \"77-2=\"], [\"97-56=\", \"22-20=\"], [\"72-48=\", \"22+57=\"], [\"9+68=\", \"91-35=\"], [\"24-1=\", \"41+18=\"], [\"77-22=\", \"35-27=\"], [\"4+36=\", \"17+20=\"], [\"24+45=\", \"29+70=\"], [\"57-12=\", \"40-35=\"], [\"80-52=\", \"59-30=\"], [\"2+58=\", \"0+79=\"], [\"95-76=\", \"82-31=\"], [\"72-3=\", \"43+52=\"], [\"4+28=\", \"6+6=\"], [\"83-26=\", \"85-36=\"], [\"16+11=\", \"9+22=\"], [\"13+22=\", \"66-18=\"], [\"59+13=\", \"67-24=\"], [\"85-70=\", \"25-12=\"], [\"40+35=\", \"81-31=\"], [\"30+41=\", \"52-5=\"], [\"79-77=\", \"7+68=\"], [\"24+17=\", \"99-40=\"], [\"2+66=\", \"20+72=\"], [\"91-69=\", \"23+46=\"], [\"15+19=\", \"45-13=\"], [\"20+75=\", \"15+77=\"], [\"41+31=\", \"78-3=\"], [\"30-25=\", \"14+53=\"], [\"74-47=\", \"44+9=\"], [\"56+22=\", \"55+24=\"], [\"2+60=\", \"13-5=\"], [\"28+55=\", \"28+26=\"], [\"93-73=\", \"41+50=\"], [\"96-69=\", \"83-6=\"], [\"28-18=\", \"71-46=\"]];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet applied = 0;\nconst mismatches = [];\n\nfor (let i = 0; i < items.length && i < REPLACEMENTS.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = items[i];\n  const currentText = para.text;\n\n  if (currentText === oldText) {\n    para.insertText(newText, \"Replace\");\n    applied++;\n  } else if (currentText === newText) {\n    // Already updated (idempotent re-run); nothing to do.\n  } else {\n    mismatches.push(i + \": expected \" + JSON.stringify(oldText) + \" got \" + JSON.stringify(currentText));\n  }\n}\n\nawait context.sync();\n\nif (mismatches.length) {\n  return \"applied=\" + applied + \" mismatches=\" + mismatches.join(\" | \");\n}\nreturn \"applied=\" + applied;\n", "ps1": "# Update the date/title paragraph and every arithmetic-problem cell in\n# the 20x5 table to match the target content. $cellValues is a FLAT list\n# of (old, new) pairs walked in row-major order (row 1 cols 1-5, row 2\n# cols 1-5, ...) alongside the Cell() traversal below -- this engine's\n# PowerShell-subset flattens nested @(...) literals, so a true 2D jagged\n# array is not reliable; a single running index sidesteps that.\n$d = $word.ActiveDocument\n\n$titleOld = '2026-01-06 Tuesday'\n$titleNew = '2026-01-07 Wednesday'\n$p1 = $d.Paragraphs.Item(1)\n$titleCurrent = $p1.Range.Text.TrimEnd([char]13)\nif ($titleCurrent -eq $titleOld) {\n    $p1.Range.Text = $titleNew\n} elseif ($titleCurrent -ne $titleNew) {\n    $p1.Range.Text = $titleNew\n}\n\n$cellValues = @(\n    @('77+6=', '69-54='),\n    @('46+47=', '40+21='),\n    @('79-0=', '25+60='),\n    @('74-36=', '61-38='),\n    @('85-84=', '29+50='),\n    @('5+52=', '77-31='),\n    @('33+8=', '25+37='),\n    @('50-10=', '3+63='),\n    @('20+65=', '71-1='),\n    @('76-29=', '52-12='),\n    @('96-45=', '89-70='),\n    @('55-3=', '54-15='),\n    @('47-37=', '33+33='),\n    @('23+5=', '21+11='),\n    @('53+5=', '41+2='),\n    @('47+2=', '76-50='),\n    @('48-10=', '18+29='),\n    @('79-51=', '9+55='),\n    @('53-42=', '32+25='),\n    @('21+15=', '59+33='),\n    @('7+60=', '49+30='),\n    @('46+11=', '42+26='),\n    @('93-77=', '51+0='),\n    @('31+52=', '40+11='),\n    @('84-1=', '53-0='),\n    @('23+60=', '96-7='),\n    @('40-4=', '42-30='),\n    @('90-63=', '86-1='),\n    @('35+32=', '72+1='),\n    @('50-20=', '73-32='),\n    @('49+7=', '88-32='),\n    @('87-66=', '68+12='),\n    @('88-51=', '82+3='),\n    @('64+20=', '6+51='),\n    @('13+68=', '70+14='),\n    @('5+46=', '93-87='),\n    @('15+38=', '53-8='),\n    @('15+83=', '28-7='),\n    @('72+21=', '68-31='),\n    @('66-18=', '2-0='),\n    @('29-0=', '14+25='),\n    @('91-64=', '89+2='),\n    @('42+8=', '85-42='),\n    @('22+56=', '53-25='),\n    @('41-9=', '66+27='),\n    @('21-18=', '31+19='),\n    @('52+35=', '33+25='),\n    @('46-7=', '30+39='),\n    @('85-21=', '47-6='),\n    @('54-5=', '47+7='),\n    @('93+3=', '8+57='),\n    @('53+8=', '20-12='),\n    @('66-30=', '29+19='),\n    @('96-59=', '72-18='),\n    @('81-41=', '38+7='),\n    @('99-37=', '77-6='),\n    @('62-40=', '18+4='),\n    @('98-6=', '56-51='),\n    @('53-1=', '33+15='),\n    @('95-63=', '30-11='),\n    @('11+65=', '17+50='),\n    @('19+21=', '16-15='),\n    @('46+43=', '19+8='),\n    @('66-19=', '19-19='),\n    @('8+59=', '77-2='),\n    @('97-56=', '22-20='),\n    @('72-48=', '22+57='),\n    @('9+68=', '91-35='),\n    @('24-1=', '41+18='),\n    @('77-22=', '35-27='),\n    @('4+36=', '17+20='),\n    @('24+45=', '29+70='),\n    @('57-12=', '40-35='),\n    @('80-52=', '59-30='),\n    @('2+58=', '0+79='),\n    @('95-76=', '82-31='),\n    @('72-3=', '43+52='),\n    @('4+28=', '6+6='),\n    @('83-26=', '85-36='),\n    @('16+11=', '9+22='),\n    @('13+22=', '66-18='),\n    @('59+13=', '67-24='),\n    @('85-70=', '25-12='),\n    @('40+35=', '81-31='),\n    @('30+41=', '52-5='),\n    @('79-77=', '7+68='),\n    @('24+17=', '99-40='),\n    @('2+66=', '20+72='),\n    @('91-69=', '23+46='),\n    @('15+19=', '45-13='),\n    @('20+75=', '15+77='),\n    @('41+31=', '78-3='),\n    @('30-25=', '14+53='),\n    @('74-47=', '44+9='),\n    @('56+22=', '55+24='),\n    @('2+60=', '13-5='),\n    @('28+55=', '28+26='),\n    @('93-73=', '41+50='),\n    @('96-69=', '83-6='),\n    @('28-18=', '71-46=')\n)\n\n$t = $d.Tables.Item(1)\n$numRows = $t.Rows.Count\n$numCols = $t.Columns.Count\n$applied = 0\n$mismatches = New-Object System.Collections.ArrayList\n$idx = 0\n\nfor ($r = 1; $r -le $numRows; $r++) {\n    for ($c = 1; $c -le $numCols; $c++) {\n        $pair = $cellValues[$idx]\n        $old = $pair[0]\n        $new = $pair[1]\n        $cellRange = $t.Cell($r, $c).Range\n        $current = $cellRange.Text\n        $currentTrimmed = $current.TrimEnd([char]13, [char]7)\n        if ($currentTrimmed -eq $old) {\n            $cellRange.Text = $new\n            $applied++\n        } elseif ($currentTrimmed -ne $new) {\n            [void]$mismatches.Add(\"R$r C$c expected '$old' got '$currentTrimmed'\")\n        }\n        $idx++\n    }\n}\n\nWrite-Output \"applied=$applied\"\nif ($mismatches.Count -gt 0) {\n    Write-Output ($mismatches -join ' | ')\n}\n"}
